$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.249.97'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.70%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.494.62'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.14%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.11'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.10'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.77%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.39%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.19'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.14%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.087.59'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.24%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.491.90'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.75'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -6.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.291.78'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.62%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '9.93'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.78'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.69'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -4.92%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '386.28'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.53%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.566'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.88%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.633.28'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '74.07'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.40%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '5.73'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.37%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.29%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.42'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.80%  '
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.23'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.16%  '
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.48'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -6.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.516.35'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.73%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +4.44%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.89%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.19'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.61%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.32%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '162.95'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0780'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.22%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.54%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '25.60'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -4.54%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '41.73'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.56%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.62%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.65'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.32%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.78%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.480.65'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.87%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.901'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.70%  '
